$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order identifiers refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1650996106457352"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961079533937"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961079533937"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961080093818"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961080733855"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961064253502.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961064413476.csv"
$ws1.Range("B4").Value = "go_stims-16509961064413476.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996106457352.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_7-165099610656935.csv"
$ws2.Range("B3").Value = "OB-16509961072813475.csv"
$ws2.Range("B4").Value = "OB-16509961070013518.csv"
$ws2.Range("B5").Value = "TB-1650996107937352.csv"
$ws2.Range("B6").Value = "OB-165099610680938.csv"
$ws2.Range("B7").Value = "TB-1650996107417388.csv"
$ws2.Range("B8").Value = "ZB-match_8-16509961065293858.csv"
$ws2.Range("B9").Value = "TB-1650996107721382.csv"
$ws2.Range("B10").Value = "ZB-match_4-16509961067213483.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961079693484.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961079533937.csv"
$ws4.Range("B4").Value = "MM_stims-16509961079933825.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961079693484.csv"
$ws4.Range("B6").Value = "MM_stims-16509961080093818.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961079933825.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650996108057349.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961080253592.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961080093818.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650996108041382.csv"
